$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Overview" (sheet1): update status/date for the remaining
# record and drop the second record (ff3eefd1-...) entirely.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-40-19 14:40:33"

# Hyperlink items can't be removed individually in this host, so clear
# the whole collection and re-add only the ones that should survive.
$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0e182a82b26f2c585edf3c12755fcb37875544fc/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.md")

# ------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): update status/handoff-datetime for the
# remaining record and drop the second record entirely.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-19 14:40:30"

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0e182a82b26f2c585edf3c12755fcb37875544fc/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0e182a82b26f2c585edf3c12755fcb37875544fc/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ae5bdf519ee9578e457e3a96cb988fd154d2001/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.zh-cn.xlf", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7206b09a67a0be9b69511bdfb5cd3405bf741c27/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2967b3e9b383fcb968ff4edff2a8e68f7bf80881/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.zh-cn.xlf", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.zh-cn.xlf")

# ------------------------------------------------------------------
# Sheet "de-de" (sheet3): same treatment as "zh-cn".
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-19 14:40:33"

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0e182a82b26f2c585edf3c12755fcb37875544fc/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0e182a82b26f2c585edf3c12755fcb37875544fc/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/db9456bebbf73ca61874a99522073d84f067c23e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.de-de.xlf", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d84247da83ca612994eacb322e772474fc0e84e3/e2e/3b87fe38-2958-47c8-8120-671ff8b6381b.md", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ca0a62edc416d392f93956b0f6f7e2679b5dd1c2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.de-de.xlf", "", "", "3b87fe38-2958-47c8-8120-671ff8b6381b.4e940690aa5be81be03e99475056b2a5c10bdbfe.de-de.xlf")
